# Applies the "add null forgiving, ??, ??= operators" edit described by the
# diff: extends the "Delegate combination" bullet with a "Delegate removal"
# sentence, appends three more bulleted notes, then a small non-bulleted
# "if (x is null) { variable = expression; }" code snippet, and finally
# relocates the _GoBack bookmark from the first bullet to the end of the
# "Variable = expression;" line.

$d = $word.ActiveDocument

# --- locate the anchor paragraph -------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Delegate combination*") {
        $target = $p
        break
    }
}

# --- 1. extend the existing bullet's text, drop its (_GoBack) bookmark ----------
$target.Range.Text = "Delegate combination: if left operand id null, return value of other operand. Delegate removal: if right operand is not sublist of left operand, return left operand, if left is null, return null, if right is null, return left."
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- helper: append a new bulleted (numId=1) paragraph after $cur ---------------
function Add-BulletAfter($cur, $text) {
    $cur.Range.InsertParagraphAfter()
    $new = $cur.Next()
    $new.Range.Text = $text
    return $new
}

# --- helper: append a new non-bulleted, code-indented paragraph after $cur ------
function Add-CodeLineAfter($cur, $text, $hasLeftIndent) {
    $cur.Range.InsertParagraphAfter()
    $new = $cur.Next()
    $new.Range.Text = $text
    $new.Range.ListFormat.RemoveNumbers()
    if ($hasLeftIndent) {
        $new.Range.ParagraphFormat.LeftIndent = 21
        $new.Range.ParagraphFormat.CharacterUnitLeftIndent = 0
    } else {
        $new.Range.ParagraphFormat.CharacterUnitLeftIndent = 0
    }
    $new.Range.ParagraphFormat.FirstLineIndent = 21
    return $new
}

# --- 2-4. three more bulleted notes ---------------------------------------------
$cur = $target
$cur = Add-BulletAfter $cur "Conditional operator: A?b:c?d?e=a?b:(c?d:e)"
$cur = Add-BulletAfter $cur "X!: null forgiving operator"
$cur = Add-BulletAfter $cur "?? returns the value of left operand if it isn’t null otherwise evaluate the right operand and returns its result. ?? doesn’t evaluate its right operand if left one is not null. ??= assign the value of right-hand operand to left one only if the left one is null. ??= can be replace is: variable ??= expression"

# --- 5-8. the small "if (null) { ... }" code block -------------------------------
$cur = Add-CodeLineAfter $cur "If (variable is null)" $false
$cur = Add-CodeLineAfter $cur "{" $false
$cur = Add-CodeLineAfter $cur "Variable = expression;" $true

# move the _GoBack bookmark to the end of this paragraph's text (before the
# paragraph mark), matching its original position relative to the text it follows
$bmRange = $d.Range($cur.Range.End - 1, $cur.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cur = Add-CodeLineAfter $cur "} " $false

Write-Output "edit applied"
